$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 125, shifting the existing rows 125:148 down to 126:149.
$ws.Rows.Item(125).Insert()

# The row that used to be 125 is now at 126; copy its contents into the
# newly-inserted (currently blank) row 125 so every column starts out
# identical to its neighbour, then patch the two cells that actually differ
# (Fecha / Volumen) for the new record.
$src = $ws.Range("A126:T126")
$dst = $ws.Range("A125:T125")
$src.Copy()
$dst.PasteSpecial()

$ws.Range("D125").Value = 44476
$ws.Range("M125").Value = 80
